$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update K7 and L7 values (0 -> 2)
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 2

# Update the active selection to L8
$ws.Range("L8").Select()
